$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SamplePrep")

# Fix typo'd sample names: rows 3 and 4 erroneously duplicated "C2"/"C2_prep"
# instead of having three distinct cold-sample-prep entries. Rename the
# whole series from C1/C2/C2(dup) to the corrected CC1/CC2/CC3 (+ _prep).
$ws.Range("A2").Value = "CC1"
$ws.Range("A3").Value = "CC2"

$ws.Range("Z2").Value = "CC1_prep"
$ws.Range("Z3").Value = "CC2_prep"
$ws.Range("Z4").Value = "CC3_prep"

$ws.Range("A4").Value = "CC3"

$ws.Activate()
$ws.Range("A4").Select()
